$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1780343333333333
$ws.Range("H2").Value = 0.534103
$ws.Range("I2").Value = 0.003649670474736916
$ws.Range("J2").Value = 0.003649670474736915
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.063701
$ws.Range("N2").Value = 0.191103
$ws.Range("O2").Value = 0.01611173663836548
$ws.Range("P2").Value = 0.01611173663836548
$ws.Range("Q2").Value = 0.01134096506766667
$ws.Range("R2").Value = 0.102068685609
$ws.Range("S2").Value = 0.0000588025295057795
$ws.Range("T2").Value = 0.00005880252950577948

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1780343333333333
$ws.Range("H3").Value = 0.534103
$ws.Range("I3").Value = 0.003649670474736916
$ws.Range("J3").Value = 0.003649670474736915
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.701496333333333
$ws.Range("N3").Value = 8.104489000000001
$ws.Range("O3").Value = 0.68328279700753
$ws.Range("P3").Value = 0.68328279700753
$ws.Range("Q3").Value = 0.4809590987074445
$ws.Range("R3").Value = 4.328631888367
$ws.Range("S3").Value = 0.00249375705013404
$ws.Range("T3").Value = 0.002493757050134039

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1780343333333333
$ws.Range("H4").Value = 0.534103
$ws.Range("I4").Value = 0.003649670474736916
$ws.Range("J4").Value = 0.003649670474736915
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.188504333333333
$ws.Range("N4").Value = 3.565513
$ws.Range("O4").Value = 0.3006054663541045
$ws.Range("P4").Value = 0.3006054663541044
$ws.Range("Q4").Value = 0.2115945766487778
$ws.Range("R4").Value = 1.904351189839
$ws.Range("S4").Value = 0.001097110895097096
$ws.Range("T4").Value = 0.001097110895097096

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.833451
$ws.Range("H5").Value = 14.500353
$ws.Range("I5").Value = 0.09908483984804967
$ws.Range("J5").Value = 0.09908483984804965
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.063701
$ws.Range("N5").Value = 0.191103
$ws.Range("O5").Value = 0.01611173663836548
$ws.Range("P5").Value = 0.01611173663836548
$ws.Range("Q5").Value = 0.307895662151
$ws.Range("R5").Value = 2.771060959359
$ws.Range("S5").Value = 0.001596428844486398
$ws.Range("T5").Value = 0.001596428844486397

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.833451
$ws.Range("H6").Value = 14.500353
$ws.Range("I6").Value = 0.09908483984804967
$ws.Range("J6").Value = 0.09908483984804965
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.701496333333333
$ws.Range("N6").Value = 8.104489000000001
$ws.Range("O6").Value = 0.68328279700753
$ws.Range("P6").Value = 0.68328279700753
$ws.Range("Q6").Value = 13.05755015384633
$ws.Range("R6").Value = 117.517951384617
$ws.Range("S6").Value = 0.06770296651241854
$ws.Range("T6").Value = 0.06770296651241853

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.833451
$ws.Range("H7").Value = 14.500353
$ws.Range("I7").Value = 0.09908483984804967
$ws.Range("J7").Value = 0.09908483984804965
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.188504333333333
$ws.Range("N7").Value = 3.565513
$ws.Range("O7").Value = 0.3006054663541045
$ws.Range("P7").Value = 0.3006054663541044
$ws.Range("Q7").Value = 5.744577458454334
$ws.Range("R7").Value = 51.70119712608901
$ws.Range("S7").Value = 0.02978544449114472
$ws.Range("T7").Value = 0.02978544449114471

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 35.63223
$ws.Range("H8").Value = 106.89669
$ws.Range("I8").Value = 0.7304540385283456
$ws.Range("J8").Value = 0.7304540385283456
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.063701
$ws.Range("N8").Value = 0.191103
$ws.Range("O8").Value = 0.01611173663836548
$ws.Range("P8").Value = 0.01611173663836548
$ws.Range("Q8").Value = 2.26980868323
$ws.Range("R8").Value = 20.42827814907
$ws.Range("S8").Value = 0.01176888309519918
$ws.Range("T8").Value = 0.01176888309519917

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 35.63223
$ws.Range("H9").Value = 106.89669
$ws.Range("I9").Value = 0.7304540385283456
$ws.Range("J9").Value = 0.7304540385283456
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.701496333333333
$ws.Range("N9").Value = 8.104489000000001
$ws.Range("O9").Value = 0.68328279700753
$ws.Range("P9").Value = 0.68328279700753
$ws.Range("Q9").Value = 96.26033869349
$ws.Range("R9").Value = 866.3430482414102
$ws.Range("S9").Value = 0.4991066785310941
$ws.Range("T9").Value = 0.4991066785310941

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 35.63223
$ws.Range("H10").Value = 106.89669
$ws.Range("I10").Value = 0.7304540385283456
$ws.Range("J10").Value = 0.7304540385283456
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.188504333333333
$ws.Range("N10").Value = 3.565513
$ws.Range("O10").Value = 0.3006054663541045
$ws.Range("P10").Value = 0.3006054663541044
$ws.Range("Q10").Value = 42.34905976133
$ws.Range("R10").Value = 381.14153785197
$ws.Range("S10").Value = 0.2195784769020523
$ws.Range("T10").Value = 0.2195784769020523

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5521946666666667
$ws.Range("H11").Value = 1.656584
$ws.Range("I11").Value = 0.011319887201011
$ws.Range("J11").Value = 0.011319887201011
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.063701
$ws.Range("N11").Value = 0.191103
$ws.Range("O11").Value = 0.01611173663836548
$ws.Range("P11").Value = 0.01611173663836548
$ws.Range("Q11").Value = 0.03517535246133333
$ws.Range("R11").Value = 0.316578172152
$ws.Range("S11").Value = 0.0001823830413586935
$ws.Range("T11").Value = 0.0001823830413586934

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5521946666666667
$ws.Range("H12").Value = 1.656584
$ws.Range("I12").Value = 0.011319887201011
$ws.Range("J12").Value = 0.011319887201011
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.701496333333333
$ws.Range("N12").Value = 8.104489000000001
$ws.Range("O12").Value = 0.68328279700753
$ws.Range("P12").Value = 0.68328279700753
$ws.Range("Q12").Value = 1.491751867286222
$ws.Range("R12").Value = 13.425766805576
$ws.Range("S12").Value = 0.007734684188516539
$ws.Range("T12").Value = 0.007734684188516538

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5521946666666667
$ws.Range("H13").Value = 1.656584
$ws.Range("I13").Value = 0.011319887201011
$ws.Range("J13").Value = 0.011319887201011
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.188504333333333
$ws.Range("N13").Value = 3.565513
$ws.Range("O13").Value = 0.3006054663541045
$ws.Range("P13").Value = 0.3006054663541044
$ws.Range("Q13").Value = 0.656285754176889
$ws.Range("R13").Value = 5.906571787592
$ws.Range("S13").Value = 0.003402819971135771
$ws.Range("T13").Value = 0.00340281997113577

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 7.585023666666667
$ws.Range("H14").Value = 22.755071
$ws.Range("I14").Value = 0.155491563947857
$ws.Range("J14").Value = 0.1554915639478569
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.063701
$ws.Range("N14").Value = 0.191103
$ws.Range("O14").Value = 0.01611173663836548
$ws.Range("P14").Value = 0.01611173663836548
$ws.Range("Q14").Value = 0.4831735925903333
$ws.Range("R14").Value = 4.348562333313
$ws.Range("S14").Value = 0.002505239127815436
$ws.Range("T14").Value = 0.002505239127815435

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 7.585023666666667
$ws.Range("H15").Value = 22.755071
$ws.Range("I15").Value = 0.155491563947857
$ws.Range("J15").Value = 0.1554915639478569
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.701496333333333
$ws.Range("N15").Value = 8.104489000000001
$ws.Range("O15").Value = 0.68328279700753
$ws.Range("P15").Value = 0.68328279700753
$ws.Range("Q15").Value = 20.49091362374656
$ws.Range("R15").Value = 184.418222613719
$ws.Range("S15").Value = 0.1062447107253669
$ws.Range("T15").Value = 0.1062447107253669

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 7.585023666666667
$ws.Range("H16").Value = 22.755071
$ws.Range("I16").Value = 0.155491563947857
$ws.Range("J16").Value = 0.1554915639478569
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.188504333333333
$ws.Range("N16").Value = 3.565513
$ws.Range("O16").Value = 0.3006054663541045
$ws.Range("P16").Value = 0.3006054663541044
$ws.Range("Q16").Value = 9.014833496269222
$ws.Range("R16").Value = 81.13350146642301
$ws.Range("S16").Value = 0.0467416140946746
$ws.Range("T16").Value = 0.04674161409467458
